# The ENTSO-E "Actual Production - Hydro Water Reservoir" 15-minute series has
# been re-pulled for a later window: every timestamp in column A moves forward
# 4 days (2025-09-18/19 -> 2025-09-22/23) and column B gets the freshly fetched
# production (MW) readings for that new window.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift every timestamp in column A (rows 2-193) forward by exactly 4 days.
for ($r = 2; $r -le 193; $r++) {
    $oldTimestamp = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 1).Value = $oldTimestamp + 4
}

# Overwrite column B ("Actual Production (MW)") with the retrained values.
$ws.Cells.Item(2, 2).Value = 392
$ws.Cells.Item(3, 2).Value = 390
$ws.Cells.Item(4, 2).Value = 391
$ws.Cells.Item(5, 2).Value = 390
$ws.Cells.Item(6, 2).Value = 392
$ws.Cells.Item(7, 2).Value = 390
$ws.Cells.Item(8, 2).Value = 391
$ws.Cells.Item(9, 2).Value = 389
$ws.Cells.Item(10, 2).Value = 392
$ws.Cells.Item(11, 2).Value = 390
$ws.Cells.Item(12, 2).Value = 391
$ws.Cells.Item(13, 2).Value = 390
$ws.Cells.Item(14, 2).Value = 390
$ws.Cells.Item(15, 2).Value = 390
$ws.Cells.Item(16, 2).Value = 391
$ws.Cells.Item(17, 2).Value = 391
$ws.Cells.Item(18, 2).Value = 391
$ws.Cells.Item(19, 2).Value = 391
$ws.Cells.Item(20, 2).Value = 390
$ws.Cells.Item(21, 2).Value = 394
$ws.Cells.Item(22, 2).Value = 403
$ws.Cells.Item(23, 2).Value = 404
$ws.Cells.Item(24, 2).Value = 412
$ws.Cells.Item(25, 2).Value = 422
$ws.Cells.Item(26, 2).Value = 562
$ws.Cells.Item(27, 2).Value = 562
$ws.Cells.Item(28, 2).Value = 569
$ws.Cells.Item(29, 2).Value = 557
$ws.Cells.Item(30, 2).Value = 578
$ws.Cells.Item(31, 2).Value = 566
$ws.Cells.Item(32, 2).Value = 530
$ws.Cells.Item(33, 2).Value = 570
$ws.Cells.Item(34, 2).Value = 404
$ws.Cells.Item(35, 2).Value = 351
$ws.Cells.Item(36, 2).Value = 358
$ws.Cells.Item(37, 2).Value = 365
$ws.Cells.Item(38, 2).Value = 266
$ws.Cells.Item(39, 2).Value = 260
$ws.Cells.Item(40, 2).Value = 252
$ws.Cells.Item(41, 2).Value = 247
$ws.Cells.Item(42, 2).Value = 186
$ws.Cells.Item(43, 2).Value = 165
$ws.Cells.Item(44, 2).Value = 158
$ws.Cells.Item(45, 2).Value = 177
$ws.Cells.Item(46, 2).Value = 156
$ws.Cells.Item(47, 2).Value = 150
$ws.Cells.Item(48, 2).Value = 170
$ws.Cells.Item(49, 2).Value = 162
$ws.Cells.Item(50, 2).Value = 173
$ws.Cells.Item(51, 2).Value = 165
$ws.Cells.Item(52, 2).Value = 157
$ws.Cells.Item(53, 2).Value = 166
$ws.Cells.Item(54, 2).Value = 172
$ws.Cells.Item(55, 2).Value = 179
$ws.Cells.Item(56, 2).Value = 186
$ws.Cells.Item(57, 2).Value = 191
$ws.Cells.Item(58, 2).Value = 156
$ws.Cells.Item(59, 2).Value = 180
$ws.Cells.Item(60, 2).Value = 195
$ws.Cells.Item(61, 2).Value = 190
$ws.Cells.Item(62, 2).Value = 301
$ws.Cells.Item(63, 2).Value = 319
$ws.Cells.Item(64, 2).Value = 322
$ws.Cells.Item(65, 2).Value = 358
$ws.Cells.Item(66, 2).Value = 453
$ws.Cells.Item(67, 2).Value = 467
$ws.Cells.Item(68, 2).Value = 471
$ws.Cells.Item(69, 2).Value = 502
$ws.Cells.Item(70, 2).Value = 751
$ws.Cells.Item(71, 2).Value = 766
$ws.Cells.Item(72, 2).Value = 767
$ws.Cells.Item(73, 2).Value = 790
$ws.Cells.Item(74, 2).Value = 846
$ws.Cells.Item(75, 2).Value = 863
$ws.Cells.Item(76, 2).Value = 863
$ws.Cells.Item(77, 2).Value = 946
$ws.Cells.Item(78, 2).Value = 907
$ws.Cells.Item(79, 2).Value = 878
$ws.Cells.Item(80, 2).Value = 870
$ws.Cells.Item(81, 2).Value = 867
$ws.Cells.Item(82, 2).Value = 848
$ws.Cells.Item(83, 2).Value = 850
$ws.Cells.Item(84, 2).Value = 851
$ws.Cells.Item(85, 2).Value = 854
$ws.Cells.Item(86, 2).Value = 742
$ws.Cells.Item(87, 2).Value = 729
$ws.Cells.Item(88, 2).Value = 730
$ws.Cells.Item(89, 2).Value = 721
$ws.Cells.Item(90, 2).Value = 443
$ws.Cells.Item(91, 2).Value = 432
$ws.Cells.Item(92, 2).Value = 432
$ws.Cells.Item(93, 2).Value = 429
$ws.Cells.Item(94, 2).Value = 415
$ws.Cells.Item(95, 2).Value = 414
$ws.Cells.Item(96, 2).Value = 413
$ws.Cells.Item(97, 2).Value = 412
$ws.Cells.Item(98, 2).Value = 415
$ws.Cells.Item(99, 2).Value = 409
$ws.Cells.Item(100, 2).Value = 410
$ws.Cells.Item(101, 2).Value = 409
$ws.Cells.Item(102, 2).Value = 419
$ws.Cells.Item(103, 2).Value = 418
$ws.Cells.Item(104, 2).Value = 416
$ws.Cells.Item(105, 2).Value = 414
$ws.Cells.Item(106, 2).Value = 404
$ws.Cells.Item(107, 2).Value = 404
$ws.Cells.Item(108, 2).Value = 404
$ws.Cells.Item(109, 2).Value = 405
$ws.Cells.Item(110, 2).Value = 408
$ws.Cells.Item(111, 2).Value = 408
$ws.Cells.Item(112, 2).Value = 408
$ws.Cells.Item(113, 2).Value = 409
$ws.Cells.Item(114, 2).Value = 412
$ws.Cells.Item(115, 2).Value = 413
$ws.Cells.Item(116, 2).Value = 417
$ws.Cells.Item(117, 2).Value = 423
$ws.Cells.Item(118, 2).Value = 439
$ws.Cells.Item(119, 2).Value = 440
$ws.Cells.Item(120, 2).Value = 442
$ws.Cells.Item(121, 2).Value = 469
$ws.Cells.Item(122, 2).Value = 752
$ws.Cells.Item(123, 2).Value = 756
$ws.Cells.Item(124, 2).Value = 770
$ws.Cells.Item(125, 2).Value = 0
$ws.Cells.Item(126, 2).Value = 0
$ws.Cells.Item(127, 2).Value = 0
$ws.Cells.Item(128, 2).Value = 0
$ws.Cells.Item(129, 2).Value = 0
$ws.Cells.Item(130, 2).Value = 0
$ws.Cells.Item(131, 2).Value = 0
$ws.Cells.Item(132, 2).Value = 0
$ws.Cells.Item(133, 2).Value = 0
$ws.Cells.Item(134, 2).Value = 0
$ws.Cells.Item(135, 2).Value = 0
$ws.Cells.Item(136, 2).Value = 0
$ws.Cells.Item(137, 2).Value = 0
$ws.Cells.Item(138, 2).Value = 0
$ws.Cells.Item(139, 2).Value = 0
$ws.Cells.Item(140, 2).Value = 0
$ws.Cells.Item(141, 2).Value = 0
$ws.Cells.Item(142, 2).Value = 0
$ws.Cells.Item(143, 2).Value = 0
$ws.Cells.Item(144, 2).Value = 0
$ws.Cells.Item(145, 2).Value = 0
$ws.Cells.Item(146, 2).Value = 0
$ws.Cells.Item(147, 2).Value = 0
$ws.Cells.Item(148, 2).Value = 0
$ws.Cells.Item(149, 2).Value = 0
$ws.Cells.Item(150, 2).Value = 0
$ws.Cells.Item(151, 2).Value = 0
$ws.Cells.Item(152, 2).Value = 0
$ws.Cells.Item(153, 2).Value = 0
$ws.Cells.Item(154, 2).Value = 0
$ws.Cells.Item(155, 2).Value = 0
$ws.Cells.Item(156, 2).Value = 0
$ws.Cells.Item(157, 2).Value = 0
$ws.Cells.Item(158, 2).Value = 0
$ws.Cells.Item(159, 2).Value = 0
$ws.Cells.Item(160, 2).Value = 0
$ws.Cells.Item(161, 2).Value = 0
$ws.Cells.Item(162, 2).Value = 0
$ws.Cells.Item(163, 2).Value = 0
$ws.Cells.Item(164, 2).Value = 0
$ws.Cells.Item(165, 2).Value = 0
$ws.Cells.Item(166, 2).Value = 0
$ws.Cells.Item(167, 2).Value = 0
$ws.Cells.Item(168, 2).Value = 0
$ws.Cells.Item(169, 2).Value = 0
$ws.Cells.Item(170, 2).Value = 0
$ws.Cells.Item(171, 2).Value = 0
$ws.Cells.Item(172, 2).Value = 0
$ws.Cells.Item(173, 2).Value = 0
$ws.Cells.Item(174, 2).Value = 0
$ws.Cells.Item(175, 2).Value = 0
$ws.Cells.Item(176, 2).Value = 0
$ws.Cells.Item(177, 2).Value = 0
$ws.Cells.Item(178, 2).Value = 0
$ws.Cells.Item(179, 2).Value = 0
$ws.Cells.Item(180, 2).Value = 0
$ws.Cells.Item(181, 2).Value = 0
$ws.Cells.Item(182, 2).Value = 0
$ws.Cells.Item(183, 2).Value = 0
$ws.Cells.Item(184, 2).Value = 0
$ws.Cells.Item(185, 2).Value = 0
$ws.Cells.Item(186, 2).Value = 0
$ws.Cells.Item(187, 2).Value = 0
$ws.Cells.Item(188, 2).Value = 0
$ws.Cells.Item(189, 2).Value = 0
$ws.Cells.Item(190, 2).Value = 0
$ws.Cells.Item(191, 2).Value = 0
$ws.Cells.Item(192, 2).Value = 0
$ws.Cells.Item(193, 2).Value = 0
